# Auto-generated edit script: "Add data for 2025-05-03"
# Updates year-to-date (2025, column L) and, in a few cases, the prior-year
# (2024, column K) totals across the Citywide Totals, By Neighborhood, and
# individual neighborhood sheets to reflect the newly added day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 2053
$ws.Cells.Item(3, 12).Value = 2068
$ws.Cells.Item(4, 11).Value = 1762
$ws.Cells.Item(4, 12).Value = 574
$ws.Cells.Item(6, 12).Value = 1860
$ws.Cells.Item(7, 11).Value = 27552
$ws.Cells.Item(7, 12).Value = 6672

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 12).Value = 51
$ws.Cells.Item(4, 12).Value = 24
$ws.Cells.Item(8, 12).Value = 423
$ws.Cells.Item(11, 12).Value = 120
$ws.Cells.Item(15, 12).Value = 44
$ws.Cells.Item(18, 12).Value = 49
$ws.Cells.Item(19, 12).Value = 190
$ws.Cells.Item(20, 12).Value = 173
$ws.Cells.Item(29, 11).Value = 1522
$ws.Cells.Item(29, 12).Value = 343
$ws.Cells.Item(31, 12).Value = 66
$ws.Cells.Item(32, 12).Value = 11
$ws.Cells.Item(33, 12).Value = 297
$ws.Cells.Item(37, 12).Value = 238
$ws.Cells.Item(42, 12).Value = 209
$ws.Cells.Item(43, 12).Value = 53
$ws.Cells.Item(44, 12).Value = 44
$ws.Cells.Item(53, 12).Value = 81
$ws.Cells.Item(54, 12).Value = 140
$ws.Cells.Item(63, 11).Value = 90
$ws.Cells.Item(63, 12).Value = 23
$ws.Cells.Item(66, 12).Value = 14
$ws.Cells.Item(67, 11).Value = 1075
$ws.Cells.Item(67, 12).Value = 241
$ws.Cells.Item(69, 12).Value = 16
$ws.Cells.Item(71, 12).Value = 19
$ws.Cells.Item(77, 12).Value = 42
$ws.Cells.Item(79, 12).Value = 184
$ws.Cells.Item(83, 12).Value = 157
$ws.Cells.Item(85, 12).Value = 353
$ws.Cells.Item(86, 12).Value = 52
$ws.Cells.Item(88, 12).Value = 91
$ws.Cells.Item(89, 12).Value = 86
$ws.Cells.Item(90, 12).Value = 64
$ws.Cells.Item(91, 12).Value = 93
$ws.Cells.Item(95, 12).Value = 99
$ws.Cells.Item(99, 12).Value = 103
$ws.Cells.Item(101, 11).Value = 27552
$ws.Cells.Item(101, 12).Value = 6672

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 12).Value = 41
$ws.Cells.Item(7, 12).Value = 120

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(6, 12).Value = 19
$ws.Cells.Item(7, 12).Value = 86

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 12).Value = 109
$ws.Cells.Item(4, 12).Value = 29
$ws.Cells.Item(6, 12).Value = 62
$ws.Cells.Item(7, 12).Value = 353

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Cells.Item(2, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 16

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 12).Value = 25
$ws.Cells.Item(4, 12).Value = 11
$ws.Cells.Item(7, 12).Value = 81

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 12).Value = 121
$ws.Cells.Item(3, 12).Value = 144
$ws.Cells.Item(4, 12).Value = 30
$ws.Cells.Item(6, 12).Value = 112
$ws.Cells.Item(7, 12).Value = 423

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 12).Value = 51
$ws.Cells.Item(7, 12).Value = 157

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 12).Value = 75
$ws.Cells.Item(6, 12).Value = 102
$ws.Cells.Item(7, 12).Value = 297

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(6, 12).Value = 21
$ws.Cells.Item(7, 12).Value = 99

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 12).Value = 76
$ws.Cells.Item(7, 12).Value = 238

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 12).Value = 44
$ws.Cells.Item(7, 12).Value = 103

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 12).Value = 22
$ws.Cells.Item(7, 12).Value = 66

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 11).Value = 299
$ws.Cells.Item(2, 12).Value = 69
$ws.Cells.Item(3, 12).Value = 79
$ws.Cells.Item(4, 11).Value = 60
$ws.Cells.Item(6, 12).Value = 66
$ws.Cells.Item(7, 11).Value = 1075
$ws.Cells.Item(7, 12).Value = 241

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 12).Value = 27
$ws.Cells.Item(7, 12).Value = 140

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 12).Value = 110
$ws.Cells.Item(4, 11).Value = 73
$ws.Cells.Item(6, 12).Value = 92
$ws.Cells.Item(7, 11).Value = 1522
$ws.Cells.Item(7, 12).Value = 343

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 12).Value = 62
$ws.Cells.Item(3, 12).Value = 59
$ws.Cells.Item(7, 12).Value = 190

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 12).Value = 13
$ws.Cells.Item(7, 12).Value = 44

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 55
$ws.Cells.Item(6, 12).Value = 70
$ws.Cells.Item(7, 12).Value = 209

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(3, 12).Value = 30
$ws.Cells.Item(7, 12).Value = 93

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 12).Value = 64
$ws.Cells.Item(7, 12).Value = 184

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 12).Value = 57
$ws.Cells.Item(7, 12).Value = 173

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(2, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 49

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(6, 12).Value = 7
$ws.Cells.Item(7, 12).Value = 44

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 14

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 51

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 12).Value = 23
$ws.Cells.Item(3, 12).Value = 30
$ws.Cells.Item(7, 12).Value = 91

$ws = $wb.Worksheets.Item('Galewood')
$ws.Cells.Item(2, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 11

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 12).Value = 30
$ws.Cells.Item(7, 12).Value = 52

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(3, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 64

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 53

$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(2, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 19

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 42

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(3, 12).Value = 9
$ws.Cells.Item(6, 12).Value = 8
$ws.Cells.Item(7, 12).Value = 24
